$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-07-31 Wednesday"; new = "2024-08-01 Thursday"},
    @{old = "664÷6="; new = "742÷5="},
    @{old = "204÷4="; new = "540÷7="},
    @{old = "181÷8="; new = "943÷7="},
    @{old = "243÷8="; new = "891÷8="},
    @{old = "860÷2="; new = "189÷2="},
    @{old = "399÷4="; new = "720÷9="},
    @{old = "882÷3="; new = "842÷9="},
    @{old = "519÷6="; new = "330÷7="},
    @{old = "411÷9="; new = "364÷2="},
    @{old = "226÷3="; new = "126÷2="},
    @{old = "616÷5="; new = "663÷2="},
    @{old = "920÷4="; new = "827÷8="},
    @{old = "657÷7="; new = "978÷4="},
    @{old = "766÷2="; new = "984÷7="},
    @{old = "131÷4="; new = "242÷6="},
    @{old = "456÷9="; new = "568÷5="},
    @{old = "549÷8="; new = "855÷7="},
    @{old = "817÷8="; new = "317÷8="},
    @{old = "227÷9="; new = "624÷5="},
    @{old = "967÷6="; new = "601÷6="},
    @{old = "556÷2="; new = "107÷9="},
    @{old = "706÷7="; new = "237÷4="},
    @{old = "912÷6="; new = "258÷9="},
    @{old = "846÷9="; new = "748÷2="},
    @{old = "332÷3="; new = "259÷3="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
